# Update "paises" (countries) COVID tracking workbook:
#  - Refresh the "last updated" timestamp string (cell A1)
#  - Insert "Bermudas" as a new, higher-ranked row (pushing Aruba / Guayana
#    Francesa / Monaco down by one row, each keeping its own numbers)
#  - Drop the old "Bermudas" row further down the table and replace it with
#    "Maldivas" carrying refreshed totals (Togo / Guinea Ecuatorial shift up)
#  - Refresh case totals for a handful of unrelated countries elsewhere in
#    the table (Corea del Sur, Panama, Nueva Zelanda, Jamaica rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "datos actualizados" timestamp (was ...03:22, now ...03:52) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Abril de 2020 a las 03:52"

# --- Isolated numeric refreshes elsewhere in the table ---

# Row 29: Corea del Sur
$ws.Cells.Item(29, 2).Value = 10694
$ws.Cells.Item(29, 3).Value = 11
$ws.Cells.Item(29, 4).Value = 8277
$ws.Cells.Item(29, 5).Value = 2179
$ws.Cells.Item(29, 6).Value = 55
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = 238

# Row 49: Panama
$ws.Cells.Item(49, 2).Value = 4821
$ws.Cells.Item(49, 3).Value = 163
$ws.Cells.Item(49, 4).Value = 231
$ws.Cells.Item(49, 5).Value = 4449
$ws.Cells.Item(49, 6).Value = 94
$ws.Cells.Item(49, 7).Value = 5
$ws.Cells.Item(49, 8).Value = 141

# Row 73: Nueva Zelanda
$ws.Cells.Item(73, 2).Value = 1451
$ws.Cells.Item(73, 3).Value = 6
$ws.Cells.Item(73, 4).Value = 1036
$ws.Cells.Item(73, 5).Value = 401
$ws.Cells.Item(73, 6).Value = 2
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = 14

# Row 125: Jamaica
$ws.Cells.Item(125, 2).Value = 233
$ws.Cells.Item(125, 3).Value = 5
$ws.Cells.Item(125, 4).Value = 27
$ws.Cells.Item(125, 5).Value = 200
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 6

# --- Rows 143-149: Bermudas inserted, Aruba/Guayana Francesa/Monaco shift
#     down a row keeping their own numbers, then Maldivas (refreshed) takes
#     over the old Bermudas slot while Togo/Guinea Ecuatorial shift up ---

# Row 143: Bermudas (new, higher-ranked)
$ws.Cells.Item(143, 1).Value = "Bermudas"
$ws.Cells.Item(143, 2).Value = 98
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 39
$ws.Cells.Item(143, 5).Value = 54
$ws.Cells.Item(143, 6).Value = 10
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 5

# Row 144: Aruba (pushed down one row, same totals as before)
$ws.Cells.Item(144, 1).Value = "Aruba"
$ws.Cells.Item(144, 2).Value = 97
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 51
$ws.Cells.Item(144, 5).Value = 44
$ws.Cells.Item(144, 6).Value = 4
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 2

# Row 145: Guayana Francesa (pushed down one row, same totals as before)
$ws.Cells.Item(145, 1).Value = "Guayana Francesa"
$ws.Cells.Item(145, 2).Value = 97
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 76
$ws.Cells.Item(145, 5).Value = 20
$ws.Cells.Item(145, 6).Value = 1
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 1

# Row 146: Monaco (pushed down one row, same totals as before)
$ws.Cells.Item(146, 1).Value = "Monaco"
$ws.Cells.Item(146, 2).Value = 94
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 26
$ws.Cells.Item(146, 5).Value = 65
$ws.Cells.Item(146, 6).Value = 2
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 3

# Row 147: Maldivas (refreshed totals, now ranked above Togo)
$ws.Cells.Item(147, 1).Value = "Maldivas"
$ws.Cells.Item(147, 2).Value = 86
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 16
$ws.Cells.Item(147, 5).Value = 70
$ws.Cells.Item(147, 6).Value = 2
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 0

# Row 148: Togo (pushed up one row, same totals as before)
$ws.Cells.Item(148, 1).Value = "Togo"
$ws.Cells.Item(148, 2).Value = 86
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 56
$ws.Cells.Item(148, 5).Value = 24
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 6

# Row 149: Guinea Ecuatorial (pushed up one row, same totals as before)
$ws.Cells.Item(149, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(149, 2).Value = 83
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 7
$ws.Cells.Item(149, 5).Value = 76
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 0
